# Auto-generated Excel COM-interop script
# Adds "Prodx_Rate(pps)" column (J) to the 5 target-layer sheets,
# fixes the "Initial" sheet selection + unit label.

$wb = $excel.ActiveWorkbook

# --- "Initial" sheet: small corrections -------------------------------
$wsInitial = $wb.Worksheets.Item("Initial")
$wsInitial.Range("C8").Value = "g"
$wsInitial.Range("C9").Select()

# --- per-sheet Prodx_Rate(pps) values (column J) ----------------------
$prodxRates = @{}
$prodxRates["0um_420MeV"] = @(734.2822571020887,734.2822571020887,1133.977771435984,155.718478661305,155.718478661305,1313.026893241666,1313.026893241666,5787.446361395772,1246.109544688027,1246.109544688027,66.736490854845,88.98198780646,935.0343027630045,377.9925904786613,377.9925904786613,1535.481862757816,779.4966818004931,88.98198780646,88.98198780646,22.245496951615,2839.465871059801,2839.465871059801,2441.578933713842,311.0752419250227,311.0752419250227,88.98198780646,200.7520456609159,22.245496951615,22.245496951615,600.4475599948113,600.4475599948113,133.47298170969,956.7372266182383,956.7372266182383,222.45496951615,22.245496951615,22.245496951615,44.49099390323,44.49099390323,44.49099390323,44.49099390323,44.49099390323,22250.38010948243)
$prodxRates["9um_333MeV"] = @(29.47980490335971,29.47980490335971,2387.321624075756,2387.321624075756,1369.092779867688,191.709160721235,191.709160721235,2893.723180697886,2893.723180697886,3798.011674665975,191.709160721235,191.709160721235,14.70373091192113,44.12927850564276,44.12927850564276,1558.993363600986,1558.993363600986,441.2927850564276,1249.726698663899,235.1150084317032,191.709160721235,191.709160721235,14.70373091192113,88.25855701128552,88.25855701128552,14.70373091192113,14712.75571109093)
$prodxRates["6um_362MeV"] = @(676.4077934881309,676.4077934881309,392.4612063821507,35.62896666234272,35.62896666234272,3237.35280840576,3237.35280840576,5479.988273446622,605.87329095862,605.87329095862,285.7551640939163,285.7551640939163,694.4935633674926,178.1448333117136,178.1448333117136,907.9056479439618,907.9056479439618,285.7551640939163,2947.980490335971,1032.697460111558,35.62896666234272,35.62896666234272,195.3263146971073,195.3263146971073,106.8868999870282,35.62896666234272,35.62896666234272,195.3263146971073,195.3263146971073,35.62896666234272,374.375436502789,374.375436502789,35.62896666234272,17.81448333117136,17.81448333117136,17808.9309998184)
$prodxRates["3um_391MeV"] = @(182.6662757815541,182.6662757815541,81.20510675833441,20.25606226488521,20.25606226488521,2116.03507588533,2116.03507588533,3472.467816837463,750.5594499935142,750.5594499935142,1095.997654689324,1095.997654689324,3562.896666234272,730.6651031262163,730.6651031262163,40.51212452977042,345.4382046958102,345.4382046958102,101.4611690232197,2857.551640939163,1338.346971072772,60.7681867946556,60.7681867946556,1540.907593721624,1540.907593721624,730.6651031262163,81.20510675833441,81.20510675833441,750.5594499935142,750.5594499935142,81.20510675833441,222.45496951615,222.45496951615,40.51212452977042,60.7681867946556,60.7681867946556,20.25606226488521,20.25606226488521,20284.99949669218)
$prodxRates["12um_303MeV"] = @(245.9664703593203,245.9664703593203,74.69422960176418,4521.442469840447,4521.442469840447,2296.892774678947,117.376646517058,117.376646517058,618.5333298741731,618.5333298741731,544.3816733687898,21.3412084576469,21.3412084576469,160.0590634323518,160.0590634323518,42.68241691529381,1419.7329355299,1419.7329355299,383.4183214424699,117.376646517058,42.68241691529381,21.3412084576469,21.3412084576469,42.68241691529381,42.68241691529381,10674.60118396679)

foreach ($sheetName in $prodxRates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $prodxRates[$sheetName]
    $lastRow = $values.Count + 1

    # match column J width to the existing data columns (A..I)
    $ws.Columns.Item(10).ColumnWidth = 15.2

    # clone formatting from column I (header + data rows) onto column J
    $ws.Range("I1:I" + $lastRow).Copy()
    $ws.Range("J1:J" + $lastRow).PasteSpecial(-4122)

    $ws.Range("J1").Value = "Prodx_Rate(pps)"
    for ($i = 0; $i -lt $values.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 10).Value = $values[$i]
    }
}

Write-Output "done"
